# Append the newest round of "Partidos" match results (2025-10-11) that were
# uploaded to the workbook, rows 471-484 of the "Partidos" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Partidos")

$fecha = 45941   # serial date for 2025-10-11

# Each entry: jugador, equipo, posicion, goles, autogoles, arquero, goles_recibidos,
#             tarjetas_amarillas, tarjetas_rojas, asistencias, Penales_Atajados
$partidos = @(
    @("Gember Marin Sarria",        "Amarillo", "Arquero",       0, 0, $true,  3, 0, 0, 0, 0),
    @("Invitado",                   "Azul",     "Arquero",       0, 0, $true,  5, 0, 0, 0, 0),
    @("Sombra",                     "Amarillo", "Arquero",       0, 0, $true,  2, 0, 0, 0, 0),
    @("Quintero",                   "Amarillo", "Defensa",       1, 0, $false, 0, 0, 0, 0, 0),
    @("Arnul David Narvaez",        "Amarillo", "Delantero",     2, 0, $false, 0, 0, 0, 0, 0),
    @("David Fernando Velasco",     "Amarillo", "Delantero",     1, 0, $false, 0, 0, 0, 1, 0),
    @("Armando Murillo",            "Amarillo", "Defensa",       1, 0, $false, 0, 0, 0, 0, 0),
    @("Edwing Yesid Castillo",      "Amarillo", "Mediocampista", 0, 0, $false, 0, 0, 0, 3, 0),
    @("Juan David Espinal",         "Amarillo", "Mediocampista", 0, 0, $false, 0, 0, 0, 1, 0),
    @("Andres Tangarife",           "Azul",     "Delantero",     1, 0, $false, 0, 0, 0, 0, 0),
    @("Carlos Fernando Valencia",   "Azul",     "Delantero",     3, 0, $false, 0, 0, 0, 0, 0),
    @("Sebastian Giraldo",          "Azul",     "Mediocampista", 1, 0, $false, 0, 0, 0, 1, 0),
    @("Juan Diego Gomez Ceballos",  "Azul",     "Mediocampista", 0, 0, $false, 0, 0, 0, 2, 0),
    @("Bryan Andres Burgos",        "Azul",     "Mediocampista", 0, 0, $false, 0, 0, 0, 2, 0)
)

$startRow = 471
$r = $startRow
foreach ($p in $partidos) {
    $ws.Cells.Item($r, 1).Value  = $fecha      # fecha
    $ws.Cells.Item($r, 2).Value  = $p[0]       # jugador
    $ws.Cells.Item($r, 3).Value  = $p[1]       # equipo
    $ws.Cells.Item($r, 4).Value  = $p[2]       # posicion
    $ws.Cells.Item($r, 5).Value  = $p[3]       # goles
    $ws.Cells.Item($r, 6).Value  = $p[4]       # autogoles
    $ws.Cells.Item($r, 7).Value  = $p[5]       # arquero
    $ws.Cells.Item($r, 8).Value  = $p[6]       # goles_recibidos
    $ws.Cells.Item($r, 9).Value  = $p[7]       # tarjetas_amarillas
    $ws.Cells.Item($r, 10).Value = $p[8]       # tarjetas_rojas
    $ws.Cells.Item($r, 11).Value = $p[9]       # asistencias
    $ws.Cells.Item($r, 12).Value = $p[10]      # Penales_Atajados
    $r++
}

# Match the author's last selection/scroll position on the sheet.
$ws.Range("F478").Select()

Write-Host "Added $($partidos.Count) rows ($startRow-$($r-1)) to Partidos."
